# Updates the cryptocurrency price/listing table to reflect the latest
# scrape (GitHub Actions run on Tue Dec 13 18:44:12 UTC 2022).
#
# Most rows only get a refreshed "Price" (column D) value. A couple of
# blocks (rows 15-26 and 42-43) also have their Coin/Link/Volume(1h)
# columns rotated because the underlying ranking shifted by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value (all values are stored as text, matching the
# inlineStr cells already present in the sheet)
$updates = @{
    "D2"  = "270.06"
    "D3"  = "22.70"
    "D4"  = "6.333"
    "D5"  = "0.06186"
    "D6"  = "3.648"
    "D7"  = "6.671"
    "D8"  = "1.380"
    "D9"  = "0.8302"
    "D11" = "0.1610"
    "D12" = "0.08306"
    "D13" = "0.03558"
    "D14" = "0.03203"

    # Row 15 (was ProBitToken -> now BitMartToken)
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D15" = "0.09329"
    "E15" = "14BitMartTokenBMX"

    # Row 16 (was BitMartToken -> now MCDex)
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D16" = "3.869"
    "E16" = "15MCDexMCB"

    # Row 17 (was MCDex -> now BitForexToken)
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D17" = "0.001644"
    "E17" = "16BitForexTokenBF"

    # Row 18 (was BitForexToken -> now CoinExToken)
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "D18" = "0.04738"
    "E18" = "17CoinExTokenCET"

    # Row 19 (was CoinExToken -> now TigerCash)
    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D19" = "0.006352"
    "E19" = "18TigerCashTCH"

    # Row 20 (was TigerCash -> now HotbitToken)
    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "D20" = "0.005674"
    "E20" = "19HotbitTokenHTB"

    # Row 21 (was HotbitToken -> now BitKan)
    "B21" = "BitKan"
    "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "D21" = "0.001076"
    "E21" = "20BitKanKAN"

    # Row 22 (was BitKan -> now NitroEx)
    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "D22" = "0.0001500"
    "E22" = "21NitroExNTX"

    # Row 23 (was NitroEx -> now LEO)
    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D23" = "3.732"
    "E23" = "22LEOLEO"

    # Row 24 (was LEO -> now BTSEToken)
    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D24" = "2.413"
    "E24" = "23BTSETokenBTSE"

    # Row 25 (was BTSEToken -> now BitpandaEcosystemToken)
    "B25" = "BitpandaEcosystemToken"
    "C25" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "D25" = "0.3309"
    "E25" = "24BitpandaEcosystemTokenBEST"

    # Row 26 (was BitpandaEcosystemToken -> now ProBitToken)
    "B26" = "ProBitToken"
    "C26" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "D26" = "0.1257"
    "E26" = "25ProBitTokenPROB"

    "D27" = "0.0002704"

    "D40" = "0.04712"
    "D41" = "0.006977"

    # Row 42 (was BKEXToken -> now CEJI)
    "B42" = "CEJI"
    "C42" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D42" = "0.003801"
    "E42" = "41CEJICEJIWorstin24h"

    # Row 43 (was CEJI -> now BKEXToken)
    "B43" = "BKEXToken"
    "C43" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D43" = "0.1160"
    "E43" = "42BKEXTokenBKK"

    "D44" = "0.01184"
    "D45" = "0.00006271"
    "D46" = "0.0009900"
    "D48" = "0.9200"
    "D49" = "0.002304"
    "D50" = "0.00002301"
    "D51" = "0.01240"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        # Column D holds numeric-looking price strings (e.g. "0.1160",
        # "0.0001500") that must stay literal text, including trailing
        # zeros, exactly as scraped - force text formatting first so
        # Excel does not reinterpret/round them as numbers.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$cellRef]
}
